$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (years 2010..2021) that will occupy rows 2..13 after the shift.
$data = @(
    @("2010年", 4.6596987942, 9547.0414195352, 39719872, 19386.387, 11209484),
    @("2011年", 5.6645192064, 12025.2947372796, 46327126.9, 21822.409, 10054231),
    @("2012年", 5.6887578821, 13374.165567635, 57070803.7, 24275.3077, 10157280),
    @("2013年", 5.614395523, 13457.8084655949, 60941735.6, 25423.9766, 11467280),
    @("2014年", 6.5246534479, 12518.8488733934, 56797668.0627024, 29602.1706551259, 12030587),
    @("2015年", 5.2490265951, 11115.761282226, 56619984.2, 26736.7925, 9643725),
    @("2016年", 4.9, 10805.4143724989, 56021121.4, 25365.25, 9579327),
    @("2017年", 4.6116253736, 9913.5742964828, 54818436.5, 25500.6, 10225555),
    @("2018年", 4.9, 11774.1, 62464586.3, 25757.6541, 10923709),
    @("2019年", 4.6, 10770.1, 58450114.4, 25117.801, 9833258),
    @("2020年", 5.1, 9781.1, 52494828.4, 27321.7713, 8166436),
    @("2021年", 4.3, 8639.1, 45635348.7, 22920.8585, 7394261)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}

# Remove the now-unused rows 14..22 (old data for years that are no
# longer represented once rows 2..13 above hold 2010..2021).
$ws.Range("A14:F22").EntireRow.Delete()
